$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells are stored as plain text in the source data (e.g. "43.049.30",
# "0.0792", "1.00"). Force text format on each D cell being updated so Excel does
# not silently reinterpret the new value as a number (which would e.g. turn "1.00" into 1).
$priceCells = 'D2','D3','D5','D6','D7','D10','D11','D13','D14','D15','D16','D17','D18','D19','D20','D21','D22','D23','D26','D27','D28','D32','D33','D34','D35','D43','D44','D46','D47','D48','D50','D51'
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.049.30'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '2.300.13'
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '300.63'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '99.51'
$ws.Range("E6").Value = '  +1.94%  '
$ws.Range("D7").Value = '0.507'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E9").Value = '  +1.74%  '
$ws.Range("D10").Value = '36.33'
$ws.Range("E10").Value = '  +7.98%  '
$ws.Range("D11").Value = '0.0792'
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").Value = '17.93'
$ws.Range("E13").Value = '  +4.35%  '
$ws.Range("D14").Value = '6.91'
$ws.Range("E14").Value = '  +1.95%  '
$ws.Range("D15").Value = '2.657.71'
$ws.Range("E15").Value = '  -0.14%  '
$ws.Range("D16").Value = '2.297.82'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").Value = '0.799'
$ws.Range("E17").Value = '  -1.78%  '
$ws.Range("D18").Value = '42.964.47'
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("D19").Value = '12.83'
$ws.Range("E19").Value = '  +9.41%  '
$ws.Range("D20").Value = '0.0₃0905'
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = '6.12'
$ws.Range("E21").Value = '  +0.79%  '
$ws.Range("D22").Value = '67.95'
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").Value = '235.76'
$ws.Range("E23").Value = '  -0.37%  '
$ws.Range("E24").Value = '  +6.88%  '
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("D26").Value = '2.44'
$ws.Range("D27").Value = '24.92'
$ws.Range("E27").Value = '  +1.94%  '
$ws.Range("D28").Value = '169.08'
$ws.Range("E28").Value = '  +0.97%  '
$ws.Range("E29").Value = '  +1.28%  '
$ws.Range("E30").Value = '  -10.27%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").Value = '5.04'
$ws.Range("E33").Value = '  +1.58%  '
$ws.Range("D34").Value = '17.65'
$ws.Range("E34").Value = '  +5.20%  '
$ws.Range("D35").Value = '4.62'
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("E37").Value = '  -1.65%  '
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("E40").Value = '  +1.10%  '
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("E42").Value = '  +3.19%  '
$ws.Range("D43").Value = '2.29'
$ws.Range("E43").Value = '  -3.28%  '
$ws.Range("D44").Value = '1.985.31'
$ws.Range("E44").Value = '  +0.26%  '
$ws.Range("E45").Value = '  +2.27%  '
$ws.Range("D46").Value = '2.91'
$ws.Range("E46").Value = '  +1.70%  '
$ws.Range("D47").Value = '17.56'
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("D48").Value = '55.75'
$ws.Range("E48").Value = '  +4.38%  '
$ws.Range("E49").Value = '  +4.12%  '
$ws.Range("D50").Value = '2.524.11'
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("D51").Value = '70.72'
$ws.Range("E51").Value = '  +0.69%  '
